$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = 0.977669497583861
$ws.Cells.Item(2, 10).Value = 0.977669497583861
$ws.Cells.Item(2, 13).Value = 1.126151
$ws.Cells.Item(2, 14).Value = 3.378452999999999
$ws.Cells.Item(2, 15).Value = 0.1133308523309815
$ws.Cells.Item(2, 16).Value = 0.1133308523309815
$ws.Cells.Item(2, 17).Value = 8.462793904044998
$ws.Cells.Item(2, 18).Value = 76.16514513640499
$ws.Cells.Item(2, 19).Value = 0.1108001174591815
$ws.Cells.Item(2, 20).Value = 0.1108001174591814

$ws.Cells.Item(3, 9).Value = 0.977669497583861
$ws.Cells.Item(3, 10).Value = 0.977669497583861
$ws.Cells.Item(3, 13).Value = 6.800141333333333
$ws.Cells.Item(3, 15).Value = 0.6843361265743261
$ws.Cells.Item(3, 16).Value = 0.6843361265743262
$ws.Cells.Item(3, 17).Value = 51.10166809102666
$ws.Cells.Item(3, 18).Value = 459.91501281924
$ws.Cells.Item(3, 19).Value = 0.6690545570464069
$ws.Cells.Item(3, 20).Value = 0.669054557046407

$ws.Cells.Item(4, 9).Value = 0.977669497583861
$ws.Cells.Item(4, 10).Value = 0.977669497583861
$ws.Cells.Item(4, 13).Value = 2.010551666666667
$ws.Cells.Item(4, 14).Value = 6.031655
$ws.Cells.Item(4, 15).Value = 0.2023330210946923
$ws.Cells.Item(4, 16).Value = 0.2023330210946923
$ws.Cells.Item(4, 17).Value = 15.10888361190833
$ws.Cells.Item(4, 18).Value = 135.979952507175
$ws.Cells.Item(4, 19).Value = 0.1978148230782726
$ws.Cells.Item(4, 20).Value = 0.1978148230782726

$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.171642
$ws.Cells.Item(5, 8).Value = 0.514926
$ws.Cells.Item(5, 9).Value = 0.02233050241613897
$ws.Cells.Item(5, 10).Value = 0.02233050241613898
$ws.Cells.Item(5, 13).Value = 1.126151
$ws.Cells.Item(5, 14).Value = 3.378452999999999
$ws.Cells.Item(5, 15).Value = 0.1133308523309815
$ws.Cells.Item(5, 16).Value = 0.1133308523309815
$ws.Cells.Item(5, 17).Value = 0.193294809942
$ws.Cells.Item(5, 18).Value = 1.739653289478
$ws.Cells.Item(5, 19).Value = 0.002530734871800072
$ws.Cells.Item(5, 20).Value = 0.002530734871800072

$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.171642
$ws.Cells.Item(6, 8).Value = 0.514926
$ws.Cells.Item(6, 9).Value = 0.02233050241613897
$ws.Cells.Item(6, 10).Value = 0.02233050241613898
$ws.Cells.Item(6, 13).Value = 6.800141333333333
$ws.Cells.Item(6, 15).Value = 0.6843361265743261
$ws.Cells.Item(6, 16).Value = 0.6843361265743262
$ws.Cells.Item(6, 17).Value = 1.167189858736
$ws.Cells.Item(6, 18).Value = 10.504708728624
$ws.Cells.Item(6, 19).Value = 0.01528156952791917
$ws.Cells.Item(6, 20).Value = 0.01528156952791918

$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.171642
$ws.Cells.Item(7, 8).Value = 0.514926
$ws.Cells.Item(7, 9).Value = 0.02233050241613897
$ws.Cells.Item(7, 10).Value = 0.02233050241613898
$ws.Cells.Item(7, 13).Value = 2.010551666666667
$ws.Cells.Item(7, 14).Value = 6.031655
$ws.Cells.Item(7, 15).Value = 0.2023330210946923
$ws.Cells.Item(7, 16).Value = 0.2023330210946923
$ws.Cells.Item(7, 17).Value = 0.34509510917
$ws.Cells.Item(7, 18).Value = 3.10585598253
$ws.Cells.Item(7, 19).Value = 0.004518198016419724
$ws.Cells.Item(7, 20).Value = 0.004518198016419724
